$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new header "Label" in column H, row 1 (same style as existing headers:
# bold font, thin border, centered horizontal, top vertical alignment)
$ws.Range("H1").Value = "Label"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# Column H labels: 0 for Control rows (2-6, 12-16), 1 for MDD rows (7-11, 17-21)
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
for ($r = 7; $r -le 11; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
for ($r = 12; $r -le 16; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}
for ($r = 17; $r -le 21; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
